# Append new "specify" feature-number combinations produced by PipelineUtils
# to the bottom of the combinations list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> Feature_num text value, written in the same order the values were
# generated (ascending evens, then the odd tail descending) so the shared
# string table lands in the same sequence as the authored workbook.
$rowOrder = @(198, 200, 202, 204, 206, 207, 205, 203, 201, 199)
$featureNum = @{
    198 = "4, 4, 2"
    199 = "5, 5, 3"
    200 = "6, 6, 4"
    201 = "7, 7, 5"
    202 = "8, 8, 6"
    203 = "9, 9, 7"
    204 = "10, 10, 8"
    205 = "11, 11, 9"
    206 = "12, 12, 10"
    207 = "13, 13, 11"
}

foreach ($r in $rowOrder) {
    $ws.Range("A$r").Value = "BCE"
    $ws.Range("B$r").Value = 2
    $ws.Range("C$r").Value = "GIN"
    $ws.Range("D$r").Value = "mean"
    $ws.Range("E$r").Value = "specify"
    $ws.Range("F$r").Value = $featureNum[$r]
}

# Match the author's final viewport/selection state.
$ws.Range("A168").Select()
$ws.Range("H201").Select()
